$wb = $excel.ActiveWorkbook

# --- Jaana sheet: add two new time-tracking entries ---
$jaana = $wb.Worksheets.Item("Jaana")

# Row 9: copy the date-number-format from the row above so the new date
# cell picks up the same style (s=11) instead of minting a new one.
$jaana.Range("A7").Copy()
$jaana.Range("A9").PasteSpecial(-4122)
$jaana.Range("A9").Value = 44978
$jaana.Range("B9").Value = 0.75
$jaana.Range("C9").Value = "Vaatimusmäärittelydok., luvut 1 ja 4"

# Row 10
$jaana.Range("A7").Copy()
$jaana.Range("A10").PasteSpecial(-4122)
$jaana.Range("A10").Value = 44979
$jaana.Range("B10").Value = 0.25
$jaana.Range("C10").Value = "Vaatimusmäärittelydok., luvut 1 ja 4; ulkoasun pohtimista"

$excel.CutCopyMode = 0

# --- Selection / active sheet bookkeeping ---
$jaana.Activate()
$jaana.Range("A11").Select()

$jarno = $wb.Worksheets.Item("Jarno")
$jarno.Range("C14").Select()

$jaana.Activate()
